# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.988.38"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.884.82"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D5").Value = "'331.41"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4580"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("D8").Value = "'0.4066"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'47.51"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'0.07973"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'0.9878"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "'21.59"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").Value = "1.892.86"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'5.893"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'7.036"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'88.20"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "'0.00001028"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "'0.06546"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "'17.42"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "29.024.82"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "'5.394"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'2.205"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").Value = "2.110.79"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'156.59"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'19.55"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "'5.382"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").Value = "'117.10"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "'0.9728"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "'0.09315"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").Value = "'3.604"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").Value = "'1.400"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "'5.263"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'0.06033"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").Value = "'0.02213"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").Value = "'8.240"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").Value = "'1.182"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'0.5752"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1814"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'10.06"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").Value = "'1.257"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.07665"
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'11.99"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5439"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'2.238"
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").Value = "'1.887"
$ws.Range("E50").Value = "  -3.99%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'110.63"
$ws.Range("E51").Value = "  -1.61%  "
